$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.770.22"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.862.77"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "  +0.95%  "
$ws.Range("D5").Value = "'333.94"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("E6").Value = "  +0.94%  "
$ws.Range("D7").Value = "'0.4698"
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("D8").Value = "'0.3901"
$ws.Range("E8").Value = "  -1.04%  "
$ws.Range("E9").Value = "  -2.57%  "
$ws.Range("D10").Value = "'0.08001"
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D11").Value = "'1.005"
$ws.Range("E11").Value = "  -2.29%  "
$ws.Range("D12").Value = "'21.60"
$ws.Range("E12").Value = "  -2.00%  "
$ws.Range("D13").Value = "1.871.27"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").Value = "'6.000"
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").Value = "'7.159"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").Value = "'1.013"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").Value = "'88.16"
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("D18").Value = "'0.06716"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").Value = "'0.00001043"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").Value = "'16.92"
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("D21").Value = "'1.015"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").Value = "27.738.52"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "'5.469"
$ws.Range("E23").Value = "  -0.86%  "
$ws.Range("D24").Value = "'10.92"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("D25").Value = "'2.323"
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("D26").Value = "2.091.47"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").Value = "'158.37"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").Value = "'19.69"
$ws.Range("E28").Value = "  -2.05%  "
$ws.Range("D29").Value = "'2.099"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").Value = "'5.411"
$ws.Range("E30").Value = "  -2.86%  "
$ws.Range("D31").Value = "'120.95"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("D32").Value = "'0.9697"
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").Value = "'0.09444"
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("E34").Value = "  +1.41%  "
$ws.Range("D35").Value = "'5.310"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").Value = "'1.336"
$ws.Range("E36").Value = "  -7.85%  "
$ws.Range("D37").Value = "'0.06050"
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("D38").Value = "'0.02218"
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("D39").Value = "'1.205"
$ws.Range("E39").Value = "  -1.95%  "
$ws.Range("D40").Value = "'8.171"
$ws.Range("E40").Value = "  +0.63%  "
$ws.Range("D41").Value = "'1.012"
$ws.Range("E41").Value = "  +0.98%  "
$ws.Range("D42").Value = "'0.5927"
$ws.Range("E42").Value = "  -1.51%  "
$ws.Range("D43").Value = "'0.1885"
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("D44").Value = "'10.19"
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("D45").Value = "'1.260"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D46").Value = "'0.5623"
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").Value = "'12.00"
$ws.Range("E47").Value = "  -1.35%  "
$ws.Range("D48").Value = "'1.919"
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("D49").Value = "'3.310"
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("D50").Value = "'0.06763"
$ws.Range("E50").Value = "  -1.88%  "
$ws.Range("D51").Value = "'113.09"
$ws.Range("E51").Value = "  -1.42%  "
